# Generate Report for Handoff
# Replace the two UUID-named files that were previously handed back with two
# new files that are now "Ready for handoff" (one per locale: zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldUuid1 = "80eb6098-346a-450e-af2f-c3a99c7b83f9"
$oldUuid2 = "ce68870e-f5bc-44bb-9eea-10b9c930e382"
$newUuid1 = "a9d38774-1895-4a92-b9d2-ea54ed7a5a70"
$newUuid2 = "ffff0067dcba-e716-43b1-84d2-08e377e05498"

$newHash = "11d518846421860cb7a79ac95b58e1c940a6b15d"

$statusText = "Ready for handoff"
$zhDate = "2016-08-26 01:02:12"
$deDate = "2016-08-26 01:02:17"
$emptyDate = "0001-01-01 00:00:00"

$zhXlf = "$newUuid1.$newHash.zh-cn.xlf"
$deXlf = "$newUuid1.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("C2").Value = $statusText
$ws.Range("G2").Value = $zhXlf
$ws.Range("H2").Value = $zhDate
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $emptyDate

$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("C3").Value = $statusText
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = $zhXlf
$ws.Range("H3").Value = $zhDate
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $emptyDate

# Drop the hyperlink styling on I2/I3 (now empty, plain cells); touching the
# font also keeps the (now-empty) cell present in the sheet rather than the
# engine dropping it outright.
$ws.Range("I2").Font.Underline = $false
$ws.Range("I3").Font.Underline = $false
$ws.Range("J2").Font.Underline = $false
$ws.Range("J3").Font.Underline = $false

# Hyperlinks: remove all (engine removes whole collection on any delete),
# then recreate only the ones that should remain (A2, A3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid1.md", "", "", "$newUuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid2.md", "", "", "$newUuid2.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("C2").Value = $statusText
$ws.Range("G2").Value = $deXlf
$ws.Range("H2").Value = $deDate
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $emptyDate

$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("C3").Value = $statusText
$ws.Range("F3").Value = "'True"
$ws.Range("G3").Value = $deXlf
$ws.Range("H3").Value = $deDate
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $emptyDate

$ws.Range("I2").Font.Underline = $false
$ws.Range("I3").Font.Underline = $false
$ws.Range("J2").Font.Underline = $false
$ws.Range("J3").Font.Underline = $false

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid1.md", "", "", "$newUuid1.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid2.md", "", "", "$newUuid2.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newUuid1.md"
$ws.Range("B2").Value = "e2e\$newUuid1.md"
$ws.Range("A3").Value = "$newUuid2.md"
$ws.Range("B3").Value = "e2e\$newUuid2.md"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid1.md", "", "", "e2e\$newUuid1.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/788e44fe37c780b98c671ce2d3e3e982affad466/e2e/$newUuid2.md", "", "", "e2e\$newUuid2.md")

$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332
